# Updates cryptocurrency price/volume data in the cryptos worksheet
# per commit "Updated cryptos list on Wed Jul 12 15:44:06 UTC 2023 with GitHub Actions"
# Source data refreshed from coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price-column values are plain decimal numbers (e.g. "244.93") that Excel
# would otherwise auto-convert to a numeric type (losing the original text
# formatting, e.g. trailing zeros / exact decimal string). The sheet stores these
# as plain text, so force a text number format on those specific cells before
# writing the value, to faithfully preserve the literal string.

$ws.Range("D2").Value = '30.620.20'
$ws.Range("E2").Value = '  +0.02%  '
$ws.Range("D3").Value = '1.891.01'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.93'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4718'
$ws.Range("E7").Value = '  -0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2917'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06501'
$ws.Range("E9").Value = '  +0.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.46'
$ws.Range("E10").Value = '  +1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07772'
$ws.Range("E11").Value = '  +0.52%  '
$ws.Range("D12").Value = '1.889.24'
$ws.Range("E12").Value = '  +0.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7392'
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.27'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.197'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '283.95'
$ws.Range("E16").Value = '  +3.70%  '
$ws.Range("D17").Value = '30.707.97'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.10'
$ws.Range("E18").Value = '  -2.05%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007510'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").Value = '2.133.72'
$ws.Range("E21").Value = '  +0.81%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.003'
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.281'
$ws.Range("E23").Value = '  +0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.260'
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.180'
$ws.Range("E25").Value = '  -0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.36'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.91'
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.911'
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.350'
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09759'
$ws.Range("E30").Value = '  -2.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.479'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.295'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.148'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04888'
$ws.Range("E34").Value = '  +2.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.131'
$ws.Range("E35").Value = '  +1.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6938'
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  -0.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01904'
$ws.Range("E38").Value = '  +2.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.841'
$ws.Range("E39").Value = '  +3.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '75.70'
$ws.Range("E40").Value = '  +2.92%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.288'
$ws.Range("E41").Value = '  +1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.011'
$ws.Range("E42").Value = '  +1.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4271'
$ws.Range("E43").Value = '  +2.17%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8280'
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.46'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.544'
$ws.Range("E47").Value = '  +2.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.37'
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.977'
$ws.Range("E49").Value = '  +0.33%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '910.71'
$ws.Range("E50").Value = '  -1.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05757'
$ws.Range("E51").Value = '  +1.97%  '
